# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.6753301551942219, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0, 1.372039145084537),
    @(0.6753301551942219, 0.04240448674262143, 0.8054896365839992, 8.660232485948974, 0, 10.18345676446982),
    @(3.230985683306322,  1.667794583268128,  0.8054896365839992, 8.660232485948974, 1, 14.36450238910742),
    @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732, 1, 5.553084769722144),
    @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732, 1, 5.553084769722144),
    @(0.3048080303191223, 1.667794583268128,  0.1575252929769615, 0.496779210170732, 0, 2.626907116734944),
    @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
